$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 4) for the "logincustomer" test user
$ws.Range("A4").Value = "logincustomer"
$ws.Range("B4").Value = "logincustomer"
$ws.Range("C4").Value = "logincustomer"
$ws.Range("D4").Value = "qwerty"
$ws.Range("E4").Value = "abcd@gmail.com"

# Add a new "Region" column (F) for every row
$ws.Range("F1").Value = "East"
$ws.Range("F2").Value = "East"
$ws.Range("F3").Value = "East"
$ws.Range("F4").Value = "East"

# Add a new "Role" column (G) for every row
$ws.Range("G1").Value = "Administrator"
$ws.Range("G2").Value = "Administrator"
$ws.Range("G3").Value = "Customer"
$ws.Range("G4").Value = "Customer"

# Select the full used range, matching the workbook's saved selection state
$ws.Range("A1:H4").Select()
